$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.960.76'
$ws.Range("E2").Value = '  +0.86%  '

# Row 3
$ws.Range("D3").Value = '1.636.50'
$ws.Range("E3").Value = '  +0.24%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.14'
$ws.Range("E5").Value = '  +0.02%  '

# Row 6
$ws.Range("E6").Value = '  +0.00%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.38%  '

# Row 8
$ws.Range("E8").Value = '  -0.63%  '

# Row 9
$ws.Range("E9").Value = '  -0.44%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.74'
$ws.Range("E10").Value = '  +0.57%  '

# Row 11
$ws.Range("E11").Value = '  +0.55%  '

# Row 12
$ws.Range("E12").Value = '  +0.64%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.648.49'
$ws.Range("E13").Value = '  +0.96%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.861.86'
$ws.Range("E14").Value = '  +0.22%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.553'
$ws.Range("E15").Value = '  -0.69%  '

# Row 16
$ws.Range("D16").Value = '0.0₃0762'
$ws.Range("E16").Value = '  -0.20%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.10'
$ws.Range("E17").Value = '  +0.76%  '

# Row 18
$ws.Range("D18").Value = '25.941.68'
$ws.Range("E18").Value = '  +0.74%  '

# Row 19
$ws.Range("E19").Value = '  -0.20%  '

# Row 20
$ws.Range("E20").Value = '  +0.21%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.42'
$ws.Range("E21").Value = '  -0.70%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.98'

# Row 23
$ws.Range("E23").Value = '  +1.15%  '

# Row 24
$ws.Range("E24").Value = '  -0.28%  '

# Row 25
$ws.Range("E25").Value = '  -1.48%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.48'
$ws.Range("E26").Value = '  +0.12%  '

# Row 27
$ws.Range("E27").Value = '  +1.42%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.88'
$ws.Range("E28").Value = '  +0.25%  '

# Row 29
$ws.Range("E29").Value = '  +0.27%  '

# Row 30
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$ws.Range("E31").Value = '  +0.29%  '

# Row 32
$ws.Range("E32").Value = '  +0.57%  '

# Row 33
$ws.Range("E33").Value = '  +0.52%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.60'
$ws.Range("E34").Value = '  +1.29%  '

# Row 35
$ws.Range("E35").Value = '  +0.49%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.910'
$ws.Range("E36").Value = '  +1.12%  '

# Row 37
$ws.Range("D37").Value = '1.149.24'
$ws.Range("E37").Value = '  +1.99%  '

# Row 38
$ws.Range("E38").Value = '  -0.01%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.50'
$ws.Range("E39").Value = '  -0.87%  '

# Row 40
$ws.Range("E40").Value = '  +0.71%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.23%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.65'
$ws.Range("E42").Value = '  +1.59%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.73'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.806'
$ws.Range("E44").Value = '  +0.20%  '

# Row 45
$ws.Range("D45").Value = '1.771.65'
$ws.Range("E45").Value = '  +0.23%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.71'
$ws.Range("E46").Value = '  +1.37%  '

# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.46'
$ws.Range("E47").Value = '  +5.98%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0512'
$ws.Range("E48").Value = '  +1.66%  '

# Row 49
$ws.Range("E49").Value = '  +0.16%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.61'
$ws.Range("E50").Value = '  +0.43%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0967'
$ws.Range("E51").Value = '  +2.94%  '
